$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the leading "Home" / "Back to Home" / "Download Word Document"
#    hyperlink paragraphs that used to sit before the "Table of Contents"
#    heading at the top of the document.
# ---------------------------------------------------------------------------
$toc = $d.Content
$found = $toc.Find.Execute("Table of Contents", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found -and $toc.Start -gt 0) {
    $leading = $d.Range(0, $toc.Start)
    if ($leading.Text -match "Home" -and $leading.Text -match "Download Word Document") {
        $leading.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2) Fix table width formatting: some tables ("Phase 1: Opening Statements"
#    and "Phase 4: Closing Reflection") were left at the auto-fit default
#    instead of the rubric's standard full-width (100%) percentage layout
#    that the other tables in the document use.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    if ($tbl.PreferredWidthType -ne 2 -or $tbl.PreferredWidth -ne 250) {
        $tbl.PreferredWidthType = 2
        $tbl.PreferredWidth = 250
    }
}

Write-Output "Done."
